$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Target values (price, size, rooms, bathrooms, district, neighborhood)
# for rows 24-30, 33-36, 38-41, 43-46, 51 reflecting the reshuffled
# house-listing data described by the commit.
$rows = @{
    24 = @(789000,   286, 6, 5, "Nervión",  "Buhaira - Huerta del Rey")
    25 = @(1565000,  256, 6, 4, "Centro",   "Arenal - Museo - Tetuán")
    26 = @(230000,   68,  2, 2, "Triana",   "López de Gomara")
    27 = @(155000,   66,  2, 1, "Santa Justa - Miraflores - Cruz Roja", "Arroyo - Santa Justa")
    28 = @(320000,   108, 3, 2, "Triana",   "Ronda de Triana-Patrocinio-Turruñuelo")
    29 = @(205000,   75,  3, 2, "La Palmera - Los Bermejales", "Bami - Pineda")
    30 = @(280000,   83,  3, 1, "Centro",   "Arenal - Museo - Tetuán")
    33 = @(280000,   80,  3, 2, "Nervión",  "Gran Plaza - Marqués de Pickman - Ciudad Jardín")
    34 = @(550000,   124, 3, 2, "Centro",   "Arenal - Museo - Tetuán")
    35 = @(330000,   106, 3, 2, "Los Remedios", "Ramón de Carranza - Madre Rafols")
    36 = @(170000,   66,  2, 1, "Santa Justa - Miraflores - Cruz Roja", "Arroyo - Santa Justa")
    38 = @(289000,   82,  3, 2, "Triana",   "Ronda de Triana-Patrocinio-Turruñuelo")
    39 = @(300000,   100, 2, 1, "Centro",   "San Vicente")
    40 = @(349000,   126, 3, 2, "Triana",   "Ronda de Triana-Patrocinio-Turruñuelo")
    41 = @(499000,   189, 5, 2, "Los Remedios", "Ramón de Carranza - Madre Rafols")
    43 = @(368000,   90,  2, 2, "Centro",   "San Vicente")
    44 = @(380000,   189, 6, 3, "Nervión",  "Nervión")
    45 = @(330000,   140, 4, 2, "Nervión",  "Nervión")
    46 = @(372000,   99,  2, 1, "Centro",   "San Vicente")
    51 = @(239900,   79,  2, 1, "Nervión",  "Luis Montoto - Santa Justa")
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
}
